$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.756.07'
$ws.Range("E2").Value = '  -0.43%  '
$ws.Range("D3").Value = '1.632.26'
$ws.Range("E3").Value = '  -0.71%  '
$ws.Range("E4").Value = '  +0.19%  '
$ws.Range("D5").Value = '214.87'
$ws.Range("E5").Value = '  -0.50%  '
$ws.Range("E6").Value = '  -1.14%  '
$ws.Range("E7").Value = '  +0.25%  '
$ws.Range("E8").Value = '  -0.87%  '
$ws.Range("D9").Value = '0.0638'
$ws.Range("E9").Value = '  -0.93%  '
$ws.Range("D10").Value = '19.48'
$ws.Range("E10").Value = '  -4.43%  '
$ws.Range("D11").Value = '0.0784'
$ws.Range("D12").Value = '1.637.89'
$ws.Range("E12").Value = '  -0.14%  '
$ws.Range("D13").Value = '4.24'
$ws.Range("E13").Value = '  -1.26%  '
$ws.Range("D14").Value = '1.856.44'
$ws.Range("E14").Value = '  -0.84%  '
$ws.Range("D15").Value = '0.551'
$ws.Range("E15").Value = '  -2.41%  '
$ws.Range("D16").Value = '0.0₃0768'
$ws.Range("E16").Value = '  +0.08%  '
$ws.Range("D17").Value = '63.21'
$ws.Range("E17").Value = '  -0.31%  '
$ws.Range("D18").Value = '25.767.10'
$ws.Range("E18").Value = '  -0.58%  '
$ws.Range("E19").Value = '  +0.29%  '
$ws.Range("E20").Value = '  +0.83%  '
$ws.Range("D21").Value = '193.66'
$ws.Range("E21").Value = '  -0.72%  '
$ws.Range("D22").Value = '9.92'
$ws.Range("E22").Value = '  -0.26%  '
$ws.Range("D23").Value = '6.21'
$ws.Range("E23").Value = '  +1.45%  '
$ws.Range("E24").Value = '  +0.38%  '
$ws.Range("E25").Value = '  -0.82%  '
$ws.Range("D26").Value = '140.29'
$ws.Range("E26").Value = '  -0.06%  '
$ws.Range("E27").Value = '  -5.00%  '
$ws.Range("D28").Value = '6.82'
$ws.Range("E28").Value = '  -0.50%  '
$ws.Range("D29").Value = '15.53'
$ws.Range("E29").Value = '  -0.12%  '
$ws.Range("E30").Value = '  -0.51%  '
$ws.Range("D31").Value = '0.0487'
$ws.Range("E31").Value = '  -1.54%  '
$ws.Range("E32").Value = '  +0.40%  '
$ws.Range("E33").Value = '  -0.19%  '
$ws.Range("D34").Value = '1.59'
$ws.Range("E34").Value = '  -0.06%  '
$ws.Range("E35").Value = '  +0.65%  '
$ws.Range("D36").Value = '0.896'
$ws.Range("E36").Value = '  -1.61%  '
$ws.Range("E37").Value = '  -0.46%  '
$ws.Range("E38").Value = '  -2.17%  '
$ws.Range("D39").Value = '1.107.12'
$ws.Range("E39").Value = '  -2.06%  '
$ws.Range("E40").Value = '  -0.25%  '
$ws.Range("E41").Value = '  +0.22%  '
$ws.Range("E42").Value = '  +0.68%  '
$ws.Range("D43").Value = '99.76'
$ws.Range("E43").Value = '  +1.42%  '
$ws.Range("E44").Value = '  -0.34%  '
$ws.Range("E45").Value = '  -0.49%  '
$ws.Range("D46").Value = '54.98'
$ws.Range("E46").Value = '  -1.44%  '
$ws.Range("E47").Value = '  -2.20%  '
$ws.Range("D48").Value = '7.66'
$ws.Range("E48").Value = '  -1.39%  '
$ws.Range("E49").Value = '  -0.06%  '
$ws.Range("E50").Value = '  +5.83%  '
$ws.Range("E51").Value = '  +0.08%  '
